# Regenerate save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# This updates column G (header "K") values for rows 2-30 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 7
    3  = 4
    4  = 5
    5  = 3
    6  = 6
    7  = 3
    8  = 6
    9  = 6
    10 = 2
    11 = 5
    12 = 6
    13 = 5
    14 = 3
    15 = 9
    16 = 6
    17 = 7
    18 = 5
    19 = 3
    20 = 2
    21 = 2
    22 = 1
    23 = 6
    24 = 7
    25 = 1
    26 = 2
    27 = 4
    28 = 4
    29 = 5
    30 = 4
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
